# add unique logout feature
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Mandatory Unique Feature" now scores points (B36)
$ws.Range("B36").Value = 3

# "Unit Tests" now scores points (B41)
$ws.Range("B41").Value = 5

# "Implements at least 20 Unit Tests" must-have is now satisfied;
# set this last so the MIN(B9:B16)=1 condition recalculates the
# Sum Points formula in B54 using the updated range values above.
$ws.Range("B16").Value = 1

$excel.Calculate()
